$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, pushing existing rows 220-257 down to 221-258.
$ws.Rows.Item(220).EntireRow.Insert()

# Populate the newly inserted row 220 with the new weekly record.
$ws.Cells.Item(220, 1).Value = 3
$ws.Cells.Item(220, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 44637
$ws.Cells.Item(220, 5).Value = 5
$ws.Cells.Item(220, 6).Value = 100112001
$ws.Cells.Item(220, 7).Value = "Berenjena"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 110
$ws.Cells.Item(220, 11).Value = 9000
$ws.Cells.Item(220, 12).Value = 9500
$ws.Cells.Item(220, 13).Value = 9273
$ws.Cells.Item(220, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 155
$ws.Cells.Item(220, 17).Value = 60
$ws.Cells.Item(220, 18).Value = "Hortaliza"
